$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (SIC) and rename it "SID".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SID"

$data = @(
    @('Echelle', 'Critère', 'Indicateurs', 'Valeur EI', 'Justification prédiction CT', 'Incertitudes CT', 'Valeur après impact/ MC CT', 'Justification prédiction LT', 'Incertitudes LT', 'Valeur après impact/ MC LT'),
    @('SITE', 'Fonctionnalité', 'Surface totale (ha) d''habitat favorable'),
    @('SITE', 'Fonctionnalité', 'Nombre de patches d''habitat favorable'),
    @('SITE', 'Fonctionnalité', 'Estimation du nombre d’individus (faune)'),
    @('SITE', 'Fonctionnalité', 'Surface (ha) de nourrissage favorable'),
    @('SITE', 'Fonctionnalité', 'Surface (ha) de reproduction favorable'),
    @('SITE', 'Fonctionnalité', 'Estimation du nombre de couple'),
    @('SITE', 'Fonctionnalité', 'Surface (ha) de chasse favorable'),
    @('SITE', 'Fonctionnalité', 'Nombre de gîtes favorables'),
    @('SITE', 'Fonctionnalité', 'Surface (ha) de reproduction favorable'),
    @('SITE', 'Fonctionnalité', 'Nombre de mâle chanteurs'),
    @('SITE', 'Fonctionnalité', 'Nombre de pontes'),
    @('SITE', 'Fonctionnalité', 'Proportion surfacique de plante(s) hôte(s)'),
    @('SITE', 'Fonctionnalité', 'Nombre de stations / pieds (selon le type de plante)'),
    @('SITE', 'Diversité espèce', 'Nombre d''espèces'),
    @('SITE', 'Diversité espèce', 'Nombre de familles'),
    @('SITE', 'Pression', ' Surface de milieu ne générant pas de perturbation'),
    @('ELARGI', 'Représentativité', 'Surface totale (ha) d''habitat favorable'),
    @('ELARGI', 'Représentativité', 'Nombre d’observations de l’espèce'),
    @('ELARGI', 'Connectivité', 'Surface totale (ha) d''habitat favorable connecté au PS'),
    @('ELARGI', 'Connectivité', 'Nombre de zones favorables connectées entre elles grâce au site'),
)


for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Range("I2").Select() | Out-Null

$ws.Columns.AutoFit() | Out-Null
